# RP2040-Eins CPL sheet: the "R2" component was removed from the board, so
# its row (row 44, designator R2) is deleted from the BOM/CPL worksheet.
# Excel shifts every row below it up by one, which is reproduced here with a
# native row delete so dependent data (dimension, shared strings, etc.) is
# kept consistent automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the entire row that holds the "R2" designator (row 44).
$ws.Rows(44).Delete()

# The named range that spans the data table needs to shrink along with it.
foreach ($n in $wb.Names) {
    if ($n.RefersTo -like "*`$E`$64*") {
        $n.RefersTo = "=Sheet1!`$A`$1:`$E`$63"
    }
}

# Leave the selection where Excel would after a row delete: the full row
# that now occupies the deleted row's former position.
$ws.Range("A44:XFD44").Select()
